$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (46074 = 2026-02-21).
# Every data row (2 through 275) needs that date bumped by one day to 46075 (2026-02-22).
for ($r = 2; $r -le 275; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
